# Archive appraisals-rework and ui-layout work directories, update changelog
# Append four new changelog rows (21-24) to the existing changelog sheet,
# matching the date/style pattern already used by the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 21; Date = 46060; Version = "0.5.6"; Category = "Feature";     Description = "Appraisals rework: replaced matrix view with DataGrid-based grid showing appraisal status by year. Added modal for scheduling, completing, and editing appraisals with notes and outcome fields." ; Author = "Claude" },
    @{ Row = 22; Date = 46060; Version = "0.5.7"; Category = "Enhancement"; Description = "Migrated roles page to ListPage+ListRow. Added search to roles, onboarding items, and employee statuses. Fixed DataTable header text inconsistency. Fixed settings page tab padding and consolidated Save buttons with SaveBar component." ; Author = "Claude" },
    @{ Row = 23; Date = 46060; Version = "0.5.8"; Category = "Feature";     Description = "GroupBy component for dynamic list grouping. Added to My Training page with Category/Status options. CollapsibleSection now supports controlled mode. ListPage has expand/collapse all toggle icon." ; Author = "Claude" },
    @{ Row = 24; Date = 46060; Version = "0.5.9"; Category = "Fix";         Description = "Training status endpoint self-service fallback: users without training_matrix.view permission can now see their own training status. Fixed My Training RecordTrainingModal not rendering when groups active." ; Author = "Claude" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.Value = $r.Date
    $cellA.NumberFormat = "yyyy-mm-dd"

    $ws.Cells.Item($rowNum, 2).Value = $r.Version
    $ws.Cells.Item($rowNum, 3).Value = $r.Category
    $ws.Cells.Item($rowNum, 4).Value = $r.Description
    $ws.Cells.Item($rowNum, 5).Value = $r.Author
}
